# Regenerate save_data: column G ("K") values recomputed (K instead of Strike#)
# and written back into the sheet (std/mean recalculated upstream; here we
# just persist the recalculated s_vals for column G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 0
    4  = 2
    5  = 1
    6  = 1
    7  = 2
    8  = 3
    9  = 1
    10 = 1
    11 = 2
    12 = 0
    14 = 1
    15 = 0
    16 = 1
    17 = 1
    18 = 2
    19 = 0
    20 = 2
    21 = 2
    23 = 2
    25 = 3
    26 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
